$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 1818.8334
$ws.Range("I8").Value = 34.142857
$ws.Range("K8").Value = 102.428571
$ws.Range("M8").Value = 36.57142899999999
$ws.Range("H88").Value = 2608.9285
$ws.Range("I88").Value = 3047
$ws.Range("K88").Value = 3047
$ws.Range("M88").Value = -2641
$ws.Range("H91").Value = 2608.9285
$ws.Range("I91").Value = 3047
$ws.Range("K91").Value = 3047
$ws.Range("M91").Value = -1643
$ws.Range("H106").Value = 14833.167
$ws.Range("I106").Value = 9666.333000000001
$ws.Range("K106").Value = 9666.333000000001
$ws.Range("M106").Value = -9035.333000000001
$ws.Range("H112").Value = 2122
$ws.Range("J112").Value = 1876.8
$ws.Range("L112").Value = 5630.4
$ws.Range("N112").Value = -7846.4
$ws.Range("H121").Value = 4367.533
$ws.Range("J121").Value = 4367.533
$ws.Range("L121").Value = 13102.599
$ws.Range("N121").Value = -16596.599
$ws.Range("H132").Value = 1501.7894
$ws.Range("I132").Value = 1455.5883
$ws.Range("K132").Value = 4366.7649
$ws.Range("M132").Value = -1836.7649
$ws.Range("H137").Value = 1746.2142
$ws.Range("J137").Value = 1598.5
$ws.Range("L137").Value = 4795.5
$ws.Range("N137").Value = -9895.5
$ws.Range("H138").Value = 3733.7292
$ws.Range("I138").Value = 2911.5
$ws.Range("J138").Value = 3950.1052
$ws.Range("K138").Value = 8734.5
$ws.Range("L138").Value = 11850.3156
$ws.Range("M138").Value = -3594.5
$ws.Range("N138").Value = -22130.3156
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1528.8158
$ws.Range("I32").Value = 1275.4
$ws.Range("K32").Value = 1275.4
$ws.Range("M32").Value = -988.4000000000001
$ws.Range("H74").Value = 4219.409
$ws.Range("I74").Value = 3399.158
$ws.Range("K74").Value = 3399.158
$ws.Range("M74").Value = -2525.158
$ws.Range("H77").Value = 4219.409
$ws.Range("I77").Value = 3399.158
$ws.Range("K77").Value = 16995.79
$ws.Range("M77").Value = -12627.79
$ws.Range("H110").Value = 1400.48
$ws.Range("I110").Value = 1360.8
$ws.Range("J110").Value = 1559.2
$ws.Range("K110").Value = 1360.8
$ws.Range("L110").Value = 1559.2
$ws.Range("M110").Value = 684.2
$ws.Range("N110").Value = -5649.2
$ws.Range("H138").Value = 99995
$ws.Range("J138").Value = 99995
$ws.Range("L138").Value = 99995
$ws.Range("N138").Value = -110275
$ws.Range("H139").Value = 87942
$ws.Range("J139").Value = 87942
$ws.Range("L139").Value = 87942
$ws.Range("N139").Value = -98222
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6507.625
$ws.Range("I20").Value = 5061.4287
$ws.Range("J20").Value = 16631
$ws.Range("K20").Value = 5061.4287
$ws.Range("L20").Value = 16631
$ws.Range("M20").Value = -4814.4287
$ws.Range("N20").Value = -17125
$ws.Range("H139").Value = 99780
$ws.Range("J139").Value = 99780
$ws.Range("L139").Value = 99780
$ws.Range("N139").Value = -110060
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 261263.44
$ws.Range("I31").Value = 748740.3
$ws.Range("J31").Value = 4696.6577
$ws.Range("K31").Value = 748740.3
$ws.Range("L31").Value = 4696.6577
$ws.Range("M31").Value = -748445.3
$ws.Range("N31").Value = -5286.6577
$ws.Range("H34").Value = 261263.44
$ws.Range("I34").Value = 748740.3
$ws.Range("J34").Value = 4696.6577
$ws.Range("K34").Value = 748740.3
$ws.Range("L34").Value = 4696.6577
$ws.Range("M34").Value = -748538.3
$ws.Range("N34").Value = -5100.6577
$ws.Range("H74").Value = 60000
$ws.Range("J74").Value = 60000
$ws.Range("L74").Value = 60000
$ws.Range("N74").Value = -61748
$ws.Range("H77").Value = 60000
$ws.Range("J77").Value = 60000
$ws.Range("L77").Value = 180000
$ws.Range("N77").Value = -188736
$ws.Range("H99").Value = 5232.778
$ws.Range("I99").Value = 3999
$ws.Range("K99").Value = 3999
$ws.Range("M99").Value = -2501
$ws.Range("H126").Value = 5232.778
$ws.Range("I126").Value = 3999
$ws.Range("K126").Value = 11997
$ws.Range("M126").Value = -9527
$ws.Range("H134").Value = 4303.973
$ws.Range("I134").Value = 2203.1875
$ws.Range("K134").Value = 6609.5625
$ws.Range("M134").Value = -4074.5625
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 47.5
$ws.Range("I10").Value = 10
$ws.Range("K10").Value = 30
$ws.Range("M10").Value = 109
$ws.Range("H98").Value = 873.75
$ws.Range("J98").Value = 1247.5
$ws.Range("L98").Value = 3742.5
$ws.Range("N98").Value = -6738.5
$ws.Range("H107").Value = 5581.3887
$ws.Range("J107").Value = 8763.637000000001
$ws.Range("L107").Value = 26290.911
$ws.Range("N107").Value = -30130.911
$ws.Range("H137").Value = 6898.346
$ws.Range("I137").Value = 1488.1538
$ws.Range("J137").Value = 12308.538
$ws.Range("K137").Value = 4464.4614
$ws.Range("L137").Value = 36925.614
$ws.Range("M137").Value = 635.5385999999999
$ws.Range("N137").Value = -47125.614
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 3740790
$ws.Range("I2").Value = 4809564.5
$ws.Range("J2").Value = 79.166664
$ws.Range("K2").Value = 4809564.5
$ws.Range("L2").Value = 79.166664
$ws.Range("M2").Value = -4809451.5
$ws.Range("N2").Value = -305.166664
$ws.Range("H102").Value = 1531.238
$ws.Range("I102").Value = 1531.238
$ws.Range("K102").Value = 1531.238
$ws.Range("M102").Value = 90.76199999999994
$ws.Range("H126").Value = 2665.9
$ws.Range("I126").Value = 2457.375
$ws.Range("K126").Value = 7372.125
$ws.Range("M126").Value = -4902.125
$ws.Range("H132").Value = 28480.777
$ws.Range("I132").Value = 20577
$ws.Range("K132").Value = 61731
$ws.Range("M132").Value = -59201
$ws.Range("H136").Value = 34762.555
$ws.Range("J136").Value = 34857.875
$ws.Range("L136").Value = 104573.625
$ws.Range("N136").Value = -109673.625
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4266.25
$ws.Range("I7").Value = 4179.1816
$ws.Range("K7").Value = 4179.1816
$ws.Range("M7").Value = -4067.1816
$ws.Range("H100").Value = 3857.1428
$ws.Range("H122").Value = 5176.231
$ws.Range("I122").Value = 4940.1816
$ws.Range("J122").Value = 6474.5
$ws.Range("K122").Value = 14820.5448
$ws.Range("L122").Value = 19423.5
$ws.Range("M122").Value = -12370.5448
$ws.Range("N122").Value = -24323.5
$ws.Range("H126").Value = 4266.25
$ws.Range("I126").Value = 4179.1816
$ws.Range("K126").Value = 12537.5448
$ws.Range("M126").Value = -10067.5448
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
$ws.Range("H126").Value = 1764.7826
$ws.Range("I126").Value = 1561.1875
$ws.Range("J126").Value = 2230.1428
$ws.Range("K126").Value = 4683.5625
$ws.Range("L126").Value = 6690.428400000001
$ws.Range("M126").Value = -2213.5625
$ws.Range("N126").Value = -11630.4284
$ws.Range("H136").Value = 2292.5833
$ws.Range("I136").Value = 2026.65
$ws.Range("K136").Value = 6079.950000000001
$ws.Range("M136").Value = -3529.950000000001
$ws.Range("H139").Value = 71357.14
$ws.Range("I139").Value = 64825
$ws.Range("J139").Value = 73970
$ws.Range("K139").Value = 64825
$ws.Range("L139").Value = 73970
$ws.Range("M139").Value = -59685
